$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 617
$ws1.Range("F4").Value = 6472
$ws1.Range("F5").Value = 739
$ws1.Range("F6").Value = 1091
$ws1.Range("F7").Value = 85
$ws1.Range("F8").Value = 550
$ws1.Range("F9").Value = 200
$ws1.Range("F10").Value = 28
$ws1.Range("F11").Value = 730
$ws1.Range("F12").Value = 1209
$ws1.Range("F13").Value = 8
$ws1.Range("F15").Value = 203
$ws1.Range("F16").Value = 454
$ws1.Range("F19").Value = 1424
$ws1.Range("F20").Value = 680
$ws1.Range("F21").Value = 396
$ws1.Range("F24").Value = 1080
$ws1.Range("F25").Value = 171
$ws1.Range("F26").Value = 2236
$ws1.Range("F28").Value = 115
$ws1.Range("F29").Value = 405
$ws1.Range("F31").Value = 3641
$ws1.Range("F33").Value = 645

# Sheet: 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 715
$ws2.Range("F22").Value = 4092
$ws2.Range("F27").Value = 201
$ws2.Range("F29").Value = 92
$ws2.Range("F31").Value = 215
$ws2.Range("F32").Value = 35
$ws2.Range("F34").Value = 50
$ws2.Range("F35").Value = 1675

# Sheet: 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1206
$ws3.Range("F7").Value = 436
$ws3.Range("F10").Value = 818

# Sheet: 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1206
$ws4.Range("F5").Value = 436
$ws4.Range("F7").Value = 818
$ws4.Range("F8").Value = 617
$ws4.Range("F9").Value = 6472
$ws4.Range("F11").Value = 739
$ws4.Range("F12").Value = 1091
$ws4.Range("F13").Value = 715
$ws4.Range("F14").Value = 85
$ws4.Range("F15").Value = 550
$ws4.Range("F16").Value = 200
$ws4.Range("F17").Value = 28
$ws4.Range("F18").Value = 730
$ws4.Range("F24").Value = 1209
$ws4.Range("F25").Value = 8
$ws4.Range("F27").Value = 203
$ws4.Range("F34").Value = 680
$ws4.Range("F35").Value = 396
$ws4.Range("F39").Value = 1080
$ws4.Range("F40").Value = 171
$ws4.Range("F41").Value = 2236
$ws4.Range("F42").Value = 35
$ws4.Range("F43").Value = 1675
$ws4.Range("F44").Value = 1675
$ws4.Range("F45").Value = 115
$ws4.Range("F46").Value = 405
$ws4.Range("F47").Value = 3641
$ws4.Range("F51").Value = 645
